# Insert a new price record before the current row 247 ("Berenjena" /
# Vega Modelo de Temuco sheet). Excel's row insert shifts rows 247-340
# down to 248-341 (and bumps the used range from R340 to R341), which
# reproduces the diff's row-level shift. The freshly inserted (now
# blank) row 247 is then repopulated with the same record values as the
# row immediately below it (its old self, now at row 248), except for
# the date (column D), which gets the new date for the inserted entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 247..340 down to 248..341.
$ws.Rows.Item(247).Insert()

# Fill the newly blank row 247 with the same data as row 248 (its
# previous content, now shifted down one row).
$ws.Range("A247:R247").Value = $ws.Range("A248:R248").Value()

# The new record's own date differs from the row it was copied from.
$ws.Range("D247").Value = 44809
